$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266, shifting all existing rows (266..393) down to (267..394)
$ws.Rows.Item(266).Insert()

# Populate the newly inserted row 266 with the new data record
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44806
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = 100112023
$ws.Range("G266").Value = "Brócoli"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 1400
$ws.Range("K266").Value = 1500
$ws.Range("L266").Value = 1600
$ws.Range("M266").Value = 1550
$ws.Range("N266").Value = "$/unidad"
$ws.Range("O266").Value = "Región Metropolitana"
$ws.Range("P266").Value = 1550
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = "Hortaliza"
